{"js": "// Update the date title and the 25 division-problem answers in the table.\n// Approach: locate the title paragraph and the single table, then overwrite\n// each relevant paragraph's text in place (preserving existing run\n// formatting) rather than inserting/removing paragraphs or table rows \u2014\n// this keeps fonts/sizes/styles untouched and matches the net effect of the\n// source diff exactly (same number of populated cells, new values only).\n\nconst body = context.document.body;\n\n// --- 1) Title paragraph: \"2026-02-04 Wednesday\" -> \"2026-02-05 Thursday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\ntitlePara.insertText(\"2026-02-05 Thursday\", Word.InsertLocation.replace);\n\n// --- 2) Table of division answers ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Only these 5 rows (0-based) actually contain data; the rows in-between\n// are empty spacer rows used for vertical spacing in the layout.\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\n// New values for each data row, left-to-right (5 columns per row).\nconst newRowValues = [\n  [\"42\u00f77=6, 0\", \"18\u00f73=6, 0\", \"90\u00f79=10, 0\", \"86\u00f77=12, 2\", \"37\u00f75=7, 2\"],\n  [\"18\u00f72=9, 0\", \"81\u00f75=16, 1\", \"22\u00f77=3, 1\", \"68\u00f75=13, 3\", \"90\u00f72=45, 0\"],\n  [\"46\u00f78=5, 6\", \"37\u00f77=5, 2\", \"80\u00f72=40, 0\", \"82\u00f76=13, 4\", \"62\u00f78=7, 6\"],\n  [\"39\u00f79=4, 3\", \"18\u00f73=6, 0\", \"40\u00f72=20, 0\", \"61\u00f76=10, 1\", \"64\u00f79=7, 1\"],\n  [\"13\u00f72=6, 1\", \"34\u00f75=6, 4\", \"47\u00f75=9, 2\", \"13\u00f74=3, 1\", \"56\u00f78=7, 0\"],\n];\n\n// Gather all the cell-body paragraphs we need to rewrite first...\nconst cellParas = [];\nfor (let r = 0; r < dataRowIndexes.length; r++) {\n  const rowIndex = dataRowIndexes[r];\n  for (let c = 0; c < 5; c++) {\n    const cell = table.getCell(rowIndex, c);\n    const cellParagraphs = cell.body.paragraphs;\n    cellParagraphs.load(\"items\");\n    cellParas.push({ cellParagraphs, value: newRowValues[r][c] });\n  }\n}\nawait context.sync();\n\n// ...then write the new text into the first paragraph of each cell.\nfor (const { cellParagraphs, value } of cellParas) {\n  cellParagraphs.items[0].insertText(value, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the date title and the 25 division-problem answers in the table.\n# Approach: set each paragraph/cell Range.Text directly in place (preserving\n# existing run formatting) rather than inserting/removing paragraphs or\n# table rows/cells \u2014 this keeps fonts/sizes/styles untouched and matches\n# the net effect of the source diff exactly (same number of populated\n# cells, new values only).\n\n$d = $word.ActiveDocument\n\n# --- 1) Title paragraph: \"2026-02-04 Wednesday\" -> \"2026-02-05 Thursday\" ---\n$d.Paragraphs.Item(1).Range.Text = \"2026-02-05 Thursday\"\n\n# --- 2) Table of division answers ---\n$t = $d.Tables.Item(1)\n\n# Only these 5 rows (1-based) actually contain data; the rows in-between\n# are empty spacer rows used for vertical spacing in the layout.\n$dataRows = @(1, 5, 9, 13, 17)\n\n# New values for each data row, left-to-right (5 columns per row).\n$newRowValues = @(\n    @(\"42\u00f77=6, 0\", \"18\u00f73=6, 0\", \"90\u00f79=10, 0\", \"86\u00f77=12, 2\", \"37\u00f75=7, 2\"),\n    @(\"18\u00f72=9, 0\", \"81\u00f75=16, 1\", \"22\u00f77=3, 1\", \"68\u00f75=13, 3\", \"90\u00f72=45, 0\"),\n    @(\"46\u00f78=5, 6\", \"37\u00f77=5, 2\", \"80\u00f72=40, 0\", \"82\u00f76=13, 4\", \"62\u00f78=7, 6\"),\n    @(\"39\u00f79=4, 3\", \"18\u00f73=6, 0\", \"40\u00f72=20, 0\", \"61\u00f76=10, 1\", \"64\u00f79=7, 1\"),\n    @(\"13\u00f72=6, 1\", \"34\u00f75=6, 4\", \"47\u00f75=9, 2\", \"13\u00f74=3, 1\", \"56\u00f78=7, 0\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $row = $dataRows[$i]\n    $values = $newRowValues[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $t.Cell($row, $c)\n        $cell.Range.Text = $values[$c - 1]\n    }\n}\n"}
